$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the base "a gagner" input cell (row 3). Columns B:D on every row are
# driven by cascading formulas off column E, so this single edit ripples all
# the way down to row 102 once Excel recalculates.
$ws.Range("E3").Value = 10

# Column B was left at width 0 (hidden); give it a real width again, matching
# column C, while it stays hidden.
$ws.Columns("B").ColumnWidth = $ws.Columns("C").ColumnWidth

# Move the selection to where the author ended up working (bottom of the
# table) and select C97.
$ws.Activate()
$ws.Range("C97").Select()
